$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.177.58'
$ws.Range('E2').Value = '  -1.83%  '

# Row 3
$ws.Range('D3').Value = '2.142.91'
$ws.Range('E3').Value = '  -3.31%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.50%  '

# Row 6
$ws.Range('E6').Value = '  -4.70%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '69.29'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.10%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.564'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.60%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.46'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -9.30%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0882'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -7.93%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.26'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -7.08%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0990'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.61%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.55'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -6.79%  '

# Row 15
$ws.Range('D15').Value = '2.463.70'
$ws.Range('E15').Value = '  -3.38%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.17'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.60%  '

# Row 17
$ws.Range('D17').Value = '2.164.29'
$ws.Range('E17').Value = '  -2.30%  '

# Row 18
$ws.Range('E18').Value = '  -7.39%  '

# Row 19
$ws.Range('D19').Value = '41.018.74'
$ws.Range('E19').Value = '  -2.01%  '

# Row 20
$ws.Range('E20').Value = '  -8.23%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '68.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.33%  '

# Row 22
$ws.Range('E22').Value = '  -8.18%  '

# Row 23
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.51'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -14.19%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '224.40'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.38%  '

# Row 25
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.91'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -7.80%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.43%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.49'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -10.17%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.28'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -12.57%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.15'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.57%  '

# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.14'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.96%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '169.87'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.33%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.50'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.91%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '30.86'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.38%  '

# Row 34
$ws.Range('E34').Value = '  -6.51%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.04'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -11.06%  '

# Row 36
$ws.Range('E36').Value = '  -5.01%  '

# Row 37
$ws.Range('E37').Value = '  -8.13%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.11'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.32%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0283'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.71%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '11.65'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -14.94%  '

# Row 41
$ws.Range('E41').Value = '  -4.33%  '

# Row 42
$ws.Range('E42').Value = '  -7.57%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '57.05'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -12.43%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.184'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.91%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.11'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.05%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0953'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.07%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '96.63'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -8.10%  '

# Row 48
$ws.Range('E48').Value = '  -4.32%  '

# Row 49
$ws.Range('E49').Value = '  -5.94%  '

# Row 50
$ws.Range('E50').Value = '  -3.20%  '

# Row 51
$ws.Range('E51').Value = '  -11.02%  '
